$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.431.15"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.515.15"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.59"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.36"
$ws.Range("E6").Value = "  +2.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.509.90"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.194"
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.24"
$ws.Range("E11").Value = "  +8.35%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  -2.17%  "
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.086.12"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.26"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "608.04"
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.515.86"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.539.14"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.42"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.02"
$ws.Range("E23").Value = "  -10.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "98.46"
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.58"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.80"
$ws.Range("E29").Value = "  +2.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.01"
$ws.Range("E30").Value = "  -2.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.00"
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.04"
$ws.Range("E32").Value = "  -4.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "635.24"
$ws.Range("E33").Value = "  +12.34%  "
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0473"
$ws.Range("E39").Value = "  +4.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.78"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0743"
$ws.Range("E43").Value = "  +5.63%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.366.53"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.309"
$ws.Range("E45").Value = "  -5.37%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "32.05"
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.44"
$ws.Range("E50").Value = "  -1.96%  "
